# MegaSena_edt.xlsx — weekly update
# - Removes the (inert) "latest draws" highlight style from the previously
#   highlighted rows (398:402) by resetting them to the Normal style.
# - Appends six new Mega-Sena draws as rows 403:408.
# - Moves the selection to the newly appended block (B403:G408), matching
#   the "most recent results" selection convention used in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clear the highlight formatting that was applied to the previous
#    "latest 5 draws" block — it is now just historical data.
$ws.Range("A398:G402").Style = "Normal"

# 2) New draw data (Concurso, Bola1..Bola6)
$newDraws = @(
    @(403, 2946, 4,  13, 17, 21, 49, 54),
    @(404, 2947, 4,  10, 15, 37, 39, 44),
    @(405, 2948, 6,  24, 37, 52, 53, 58),
    @(406, 2949, 4,  6,  11, 38, 49, 54),
    @(407, 2950, 21, 23, 42, 49, 50, 60),
    @(408, 2951, 5,  8,  30, 31, 37, 45)
)

foreach ($row in $newDraws) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# 3) Make sure the newly appended rows use the plain/default style (no
#    leftover formatting from adjacent cells).
$ws.Range("A403:G408").Style = "Normal"

# 4) Update selection to the newly added block (keeps the "latest results"
#    selection convention seen in this workbook).
[void]$ws.Range("B403:G408").Select()
